# CDS test cases updates
# Replace the old "StatQuery" (column C) text on rows 2-4 (ParticipantsTab,
# SamplesTab, FilesTab) with the new CALL{} based stat query, and grow the
# row heights to fit the much longer text (Excel's max row height 409.5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStatQuery = @'
CALL{
        MATCH (p:participant)-->(s:study)
        OPTIONAL MATCH (samp:sample)-->(p)
        OPTIONAL MATCH (samp)<--(f:file)
        OPTIONAL MATCH (f)<--(g:genomic_info)
        OPTIONAL MATCH (p)<--(diag:diagnosis)
        WITH s, p, samp, f, g, diag
        WHERE f.file_type IN ["JSON"]
        RETURN 
            count(distinct p) AS num_participants
    }
    WITH num_participants
    CALL {
        MATCH (samp:sample)-->(p:participant)-->(s)
        OPTIONAL MATCH (samp)<--(f:file)
        OPTIONAL MATCH (p)<--(diag:diagnosis)
        OPTIONAL MATCH (f)<--(g:genomic_info)
        OPTIONAL MATCH (p)<--(diag:diagnosis)
        WITH s, p, samp, f, g, diag
        WHERE f.file_type IN ["JSON"]
        RETURN 
            count(distinct samp) AS num_samples
    }
    WITH num_participants, num_samples
    CALL {
        MATCH (f:file)-->(s:study)
        OPTIONAL MATCH (f)<--(g:genomic_info)
        OPTIONAL MATCH (samp:sample)<--(f)
        OPTIONAL MATCH (p:participant)<--(samp)
        OPTIONAL MATCH (p)<--(diag:diagnosis)
        WITH s, p, samp, f, g, diag
        WHERE f.file_type IN ["JSON"]
        RETURN 
            count(distinct s) AS num_studies,
            count(distinct f) AS num_files
    }
    RETURN 
        num_studies AS Studies,
        num_participants AS Participants,
        num_samples AS Samples,
        num_files AS `Files`
'@

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5

$ws.Range("C5").Select() | Out-Null
